# Word COM-interop script implementing the Profile.docx edit described by the diff.
# Strategy: work from the END of the affected region toward the START, using
# live $d.Paragraphs.Item(N) look-ups only (paragraph object references do not
# stay "pinned" across structural edits in this runtime, so we must always
# re-resolve by index, and process bottom-up so already-visited indices don't
# shift under us).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# [71] '(C, CPrevStatement, CMapping, CNextStatement);'
#      -> '(C, CPrevResource, CMappingResource, CNextResource);'
#      then append a brand-new paragraph after it.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "(C, CPrevStatement, CMapping, CNextStatement);", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "(C, CPrevResource, CMappingResource, CNextResource);", 2) | Out-Null

$p = $d.Paragraphs.Item(71)
$p.Range.InsertParagraphAfter()
$d.Paragraphs.Item(72).Range.Text = "Inferred / Stated. Entailment (Pattern Matching)"

# ---------------------------------------------------------------------------
# [70] 'Transforms Order:' -> 'Transforms Order (Functional Context Roles):'
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Transforms Order:", $true, $false, $false, $false, $false, $true, 1,
    $false, "Transforms Order (Functional Context Roles):", 2) | Out-Null

# ---------------------------------------------------------------------------
# [69] 'Statement: Transform (I) Quad Contexts.' -> DELETE
# [68] 'Statement: Mapping (C) Kinds'            -> DELETE
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(69).Range.Delete()
$d.Paragraphs.Item(68).Range.Delete()

# ---------------------------------------------------------------------------
# [67] 'Statement: Statement (D) SPOs' -> 'Mapping Match / Apply Transform.'
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Statement: Statement (D) SPOs", $true, $false, $false, $false, $false,
    $true, 1, $false, "Mapping Match / Apply Transform.", 2) | Out-Null

# ---------------------------------------------------------------------------
# New empty paragraph right before 'Functions:' (after 'Quad Contexts. Stream:
# Transforms (Interactions).').
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(64).Range.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# [63] 'Kinds. Stream: Mappings (Roles).' -> 'Kinds. Stream: Mappings (Resource
#      Occurrences).' then 7 brand-new paragraphs appended after it.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Kinds. Stream: Mappings (Roles).", $true, $false, $false, $false, $false,
    $true, 1, $false, "Kinds. Stream: Mappings (Resource Occurrences).", 2) | Out-Null

$newAfterKinds = @(
    "Resource Monad : Kinds (SK, S, Attribute, Value) Kind Relative Resource Roles in Contexts.",
    "Statement: Statement (D) Resource CSPOs.",
    "Resource Monad : Statements (Resource Quads). CSPO Resource Context Roles.",
    "Statement: Mapping (C) Resource Kinds.",
    "Resource Monad : Mappings (Resource Quads). Kinds Resource  Context Roles.",
    "Statement: Transform (I) Resource Quad Contexts.",
    "Resource Monad : Transform (Resource Mappings). Resource  Context Roles."
)
$anchor = $d.Paragraphs.Item(63)
foreach ($txt in $newAfterKinds) {
    $anchor.Range.InsertParagraphAfter()
    $anchor = $anchor.Next()
    $anchor.Range.Text = $txt
}

# ---------------------------------------------------------------------------
# [62] 'SPOs. Stream: Statements (Occurrences).' -> 'SPOs. Stream: Statements
#      (Context Occurrences).' then 1 new paragraph appended after it.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "SPOs. Stream: Statements (Occurrences).", $true, $false, $false, $false,
    $false, $true, 1, $false, "SPOs. Stream: Statements (Context Occurrences).",
    2) | Out-Null

$p = $d.Paragraphs.Item(62)
$p.Range.InsertParagraphAfter()
$d.Paragraphs.Item(63).Range.Text = "Resource Monad : CSPOs (URNs Resource Roles in Contexts)"

# ---------------------------------------------------------------------------
# [61] 'Resource' (whole-paragraph text, set directly -- plain "Resource" is
#      too common a substring elsewhere to use Find/Replace safely).
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(61).Range.Text = "Resource Monad Types:"

# ---------------------------------------------------------------------------
# [60] 'Monads:' -> becomes an empty paragraph (run survives, text removed).
#      Delete the paragraph and re-insert a fresh (empty) one after [59] so
#      the resulting run keeps its rPr but drops <w:t> entirely, matching the
#      target shape.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(60).Range.Delete()
$d.Paragraphs.Item(59).Range.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# [55] 'Pattern Matching: ... Signature Function Verticle.' -> '... Monad
#      Type / Instance Wrapper Verticle.'
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Yields corresponding Signature Function Verticle.", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "Yields corresponding Monad Type / Instance Wrapper Verticle.", 2) | Out-Null

# ---------------------------------------------------------------------------
# [53],[52],[51],[50] -> DELETE (old "Case"/"Reactive" block, fully replaced)
# [49] 'Switch Case Patterns...' -> DELETE (text unchanged, but its pPr carries
#      an unwanted <w:rPr><w:u w:val="none"/></w:rPr> we cannot strip in place;
#      delete + recreate from the clean [48] anchor instead).
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(53).Range.Delete()
$d.Paragraphs.Item(52).Range.Delete()
$d.Paragraphs.Item(51).Range.Delete()
$d.Paragraphs.Item(50).Range.Delete()
$d.Paragraphs.Item(49).Range.Delete()

# ---------------------------------------------------------------------------
# [48] 'Case Classes...' stays put (clean formatting, no <w:u val="none">).
#      Append the full replacement block after it: the (reinstated) "Switch
#      Case Patterns" paragraph plus the 7 new "Verticles:" / "Reactive:" /
#      "CSPO Cases:" / etc. paragraphs.
# ---------------------------------------------------------------------------
$newAfterCaseClasses = @(
    "Switch Case Patterns: CSPOs, Kinds, Statements Layers Networks. Aggregation, Alignment, Activation (Map Reduce).",
    "Verticles: URNs (CSPOs, reified Kinds, reified Statements) Monads. Encoding: extract Resources / Roles in URNs Functional Contexts Roles.",
    "Reactive / Event Driven: Verticles DIDs (Distributed IDs) distributed Resource / Applicable events logs. Rx Facade. Resource URNs Verticle Resolution, Transforms ordered Mappings Statements.",
    "CSPO Cases: one for each Statement CSPO destructuring case (for each CSPOs Data Aggregation). SCase, PCase, OCase yielding corresponding Verticle Monads.",
    "Kinds Cases: One for each Kind Type Data / Mappings Agreggation. SKCase, PKCase, OKCase yielding corresponding Verticle Monads.",
    "Statement Cases: CSPOs / Kinds Aggregation Function (Map Reduce). Entailments.",
    "Mapping Cases: Statements Aggregation Function (Map Reduce). Entailments.",
    "Transforms Cases: Mappings Aggregation Function (Map Reduce). Entailments."
)
$anchor = $d.Paragraphs.Item(48)
foreach ($txt in $newAfterCaseClasses) {
    $anchor.Range.InsertParagraphAfter()
    $anchor = $anchor.Next()
    $anchor.Range.Text = $txt
}

Write-Output "done"
